# Fix capitalization/wording in the SmartArt graphic on the "The test: the
# firewall" slide (slide 9):
#   "the Firewall is not perfect"                  -> "The firewall is not perfect"
#   "Not having a firewall Puts others at risk"     -> "Not having a firewall puts others at risk"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasSmartArt) {
            $nodes = $shape.SmartArt.AllNodes
            for ($k = 1; $k -le $nodes.Count; $k++) {
                $tr = $nodes.Item($k).TextFrame2.TextRange
                if ($tr.Text -ceq "the Firewall is not perfect") {
                    $tr.Text = "The firewall is not perfect"
                }
                elseif ($tr.Text -ceq "Not having a firewall Puts others at risk") {
                    $tr.Text = "Not having a firewall puts others at risk"
                }
            }
        }
    }
}
